$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("new_results")
$ws2 = $wb.Worksheets.Item("10-fold cross validation")

# --- new_results sheet: remove the stray B12/B13 formulas (=H8, =N8) ---
$ws1.Range("B12").ClearContents()
$ws1.Range("B13").ClearContents()

# --- chart1 (Grafiek 1 on new_results): reset the value axis to auto scaling ---
$co1 = $ws1.ChartObjects(1)
$chart1 = $co1.Chart
$valAxis1 = $chart1.Axes(2)
$valAxis1.MinimumScaleIsAuto = $true
$valAxis1.MaximumScaleIsAuto = $true

# --- chart1: move/resize the chart's anchor on the sheet ---
$co1.Left = 133.48828125
$co1.Top = 234.75
$co1.Width = 1060.3505859375
$co1.Height = 394.5

# --- new_results sheet view: scroll position + selection ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("M15").Select()

# --- 10-fold cross validation sheet view: selection ---
$ws2.Activate()
$ws2.Range("B3").Select()

$ws1.Activate()
